# "Changes of OrderCreation Time"
# B2:B11 on Sheet1 hold job/order numbers (stored as text, even though they
# look numeric) that need to be bumped to a newer batch of values.
#
# The leading "'" forces Excel to keep the entry as text (shared string)
# instead of re-interpreting the all-digit value as a number; re-applying
# the "Normal" style afterwards clears the quote-prefix formatting flag that
# the text entry leaves behind, so the cells end up with the same (default)
# style they started with.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "'32376236"
$ws.Range("B3").Value = "'32376237"
$ws.Range("B4").Value = "'32376238"
$ws.Range("B5").Value = "'32376239"
$ws.Range("B6").Value = "'32376241"
$ws.Range("B7").Value = "'32376242"
$ws.Range("B8").Value = "'32376243"
$ws.Range("B9").Value = "'32376244"
$ws.Range("B10").Value = "'32376245"
$ws.Range("B11").Value = "'32376246"

$ws.Range("B2:B11").Style = "Normal"
